# Apply the "text" category addition to the hidden '#system' sheet that backs
# the MacroLibrary workbook's defined names:
#   - base:  new function outputToCloud(resource)
#   - text (new category): new function spellCheck(var,profile,text)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------------
# 1) Column A ("target" - the list of category names): insert a new row at
#    A25 ("text"), pushing the existing entries (web..xml) down by one row.
# ---------------------------------------------------------------------------
for ($r = 30; $r -ge 25; $r--) {
    $srcVal = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r + 1, 1).Value = $srcVal
}
$ws.Cells.Item(25, 1).Value = "text"

# ---------------------------------------------------------------------------
# 2) Column E ("base" function list): insert a new row at E21
#    ("outputToCloud(resource)"), pushing the remaining entries down by
#    one row.
# ---------------------------------------------------------------------------
for ($r = 37; $r -ge 21; $r--) {
    $srcVal = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r + 1, 5).Value = $srcVal
}
$ws.Cells.Item(21, 5).Value = "outputToCloud(resource)"

# ---------------------------------------------------------------------------
# 3) Insert a new column before Y: columns Y..AD (25..30) shift right to
#    Z..AE (26..31). Column Y becomes the new "text" column.
# ---------------------------------------------------------------------------
for ($c = 30; $c -ge 25; $c--) {
    for ($r = 1; $r -le 129; $r++) {
        $srcVal = $ws.Cells.Item($r, $c).Value2
        if ($srcVal -ne $null -and $srcVal -ne "") {
            $ws.Cells.Item($r, $c + 1).Value = $srcVal
        } else {
            $ws.Cells.Item($r, $c + 1).ClearContents()
        }
    }
}

# Clear out the old "web" data that is still sitting in column Y (it has
# already been copied over to column Z above), then populate the new
# single-entry "text" column.
for ($r = 3; $r -le 129; $r++) {
    $ws.Cells.Item($r, 25).ClearContents()
}
$ws.Cells.Item(1, 25).Value = "text"
$ws.Cells.Item(2, 25).Value = "spellCheck(var,profile,text)"

# ---------------------------------------------------------------------------
# 4) Update the workbook-level defined names so they keep pointing at the
#    right ranges after the row/column insertions above.
# ---------------------------------------------------------------------------
$wb.Names.Item("base").RefersTo = "='#system'!`$E`$2:`$E`$39"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$31"
$wb.Names.Item("web").RefersTo = "='#system'!`$Z`$2:`$Z`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AC`$2:`$AC`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AD`$2:`$AD`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AE`$2:`$AE`$27"
$wb.Names.Add("text", "='#system'!`$Y`$2:`$Y`$2")

# ---------------------------------------------------------------------------
# 5) The sheet's used range grows by one column (Y..AD -> Z..AE). Nudge the
#    worksheet's extent one column past the new last column (AE) so the
#    recorded dimension keeps the same "one column wider than the data"
#    convention already present in the source file, without touching any
#    cell's font/fill formatting.
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 32).Style = $ws.Cells.Item(1, 1).Style

Write-Host "Applied 'text' category (outputToCloud/spellCheck) edit."
